$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date formatting (style used in column B) from the last existing row
$ws.Range("B222").Copy() | Out-Null
$ws.Range("B223:B230").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 223
$ws.Cells.Item(223, 1).Value = "Entrainement"
$ws.Cells.Item(223, 2).Value = 45874
$ws.Cells.Item(223, 3).Value = "Global"
$ws.Cells.Item(223, 5).Value = "Romain Thunet"
$ws.Cells.Item(223, 6).Value = "center back"
$ws.Cells.Item(223, 7).Value = "01:26:21"
$ws.Cells.Item(223, 8).Value = 4.06
$ws.Cells.Item(223, 9).Value = 0
$ws.Cells.Item(223, 10).Value = 4.06
$ws.Cells.Item(223, 11).Value = 0
$ws.Cells.Item(223, 12).Value = 0
$ws.Cells.Item(223, 13).Value = 0
$ws.Cells.Item(223, 14).Value = 0
$ws.Cells.Item(223, 15).Value = 0
$ws.Cells.Item(223, 16).Value = 2.36
$ws.Cells.Item(223, 17).Value = 8.88
$ws.Cells.Item(223, 18).Value = 0
$ws.Cells.Item(223, 19).Value = 0
$ws.Cells.Item(223, 20).Value = 0
$ws.Cells.Item(223, 21).Value = 0
$ws.Cells.Item(223, 22).Value = 0

# Row 224
$ws.Cells.Item(224, 1).Value = "Entrainement"
$ws.Cells.Item(224, 2).Value = 45874
$ws.Cells.Item(224, 3).Value = "Global"
$ws.Cells.Item(224, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(224, 6).Value = "left forward"
$ws.Cells.Item(224, 7).Value = "01:28:09"
$ws.Cells.Item(224, 8).Value = 7.15
$ws.Cells.Item(224, 9).Value = 1.66
$ws.Cells.Item(224, 10).Value = 5.47
$ws.Cells.Item(224, 11).Value = 0.77
$ws.Cells.Item(224, 12).Value = 0.9
$ws.Cells.Item(224, 13).Value = 0.01
$ws.Cells.Item(224, 14).Value = 0
$ws.Cells.Item(224, 15).Value = 1
$ws.Cells.Item(224, 16).Value = 4.8
$ws.Cells.Item(224, 17).Value = 26.86
$ws.Cells.Item(224, 18).Value = 4.46
$ws.Cells.Item(224, 19).Value = 45
$ws.Cells.Item(224, 20).Value = 6
$ws.Cells.Item(224, 21).Value = 21
$ws.Cells.Item(224, 22).Value = 3

# Row 225
$ws.Cells.Item(225, 1).Value = "Entrainement"
$ws.Cells.Item(225, 2).Value = 45874
$ws.Cells.Item(225, 3).Value = "Global"
$ws.Cells.Item(225, 5).Value = "Jeremie Laurent"
$ws.Cells.Item(225, 6).Value = "left forward"
$ws.Cells.Item(225, 7).Value = "01:26:10"
$ws.Cells.Item(225, 8).Value = 7.73
$ws.Cells.Item(225, 9).Value = 1.6
$ws.Cells.Item(225, 10).Value = 6.12
$ws.Cells.Item(225, 11).Value = 1.27
$ws.Cells.Item(225, 12).Value = 0.33
$ws.Cells.Item(225, 13).Value = 0.01
$ws.Cells.Item(225, 14).Value = 0
$ws.Cells.Item(225, 15).Value = 1
$ws.Cells.Item(225, 16).Value = 5.33
$ws.Cells.Item(225, 17).Value = 27.04
$ws.Cells.Item(225, 18).Value = 5.14
$ws.Cells.Item(225, 19).Value = 61
$ws.Cells.Item(225, 20).Value = 12
$ws.Cells.Item(225, 21).Value = 21
$ws.Cells.Item(225, 22).Value = 6

# Row 226
$ws.Cells.Item(226, 1).Value = "Entrainement"
$ws.Cells.Item(226, 2).Value = 45874
$ws.Cells.Item(226, 3).Value = "Global"
$ws.Cells.Item(226, 5).Value = "Mattheo Haon"
$ws.Cells.Item(226, 6).Value = "right back"
$ws.Cells.Item(226, 7).Value = "01:28:36"
$ws.Cells.Item(226, 8).Value = 7.05
$ws.Cells.Item(226, 9).Value = 1.35
$ws.Cells.Item(226, 10).Value = 5.69
$ws.Cells.Item(226, 11).Value = 1
$ws.Cells.Item(226, 12).Value = 0.35
$ws.Cells.Item(226, 13).Value = 0.01
$ws.Cells.Item(226, 14).Value = 0
$ws.Cells.Item(226, 15).Value = 1
$ws.Cells.Item(226, 16).Value = 4.73
$ws.Cells.Item(226, 17).Value = 26.54
$ws.Cells.Item(226, 18).Value = 4.57
$ws.Cells.Item(226, 19).Value = 38
$ws.Cells.Item(226, 20).Value = 2
$ws.Cells.Item(226, 21).Value = 17
$ws.Cells.Item(226, 22).Value = 5

# Row 227
$ws.Cells.Item(227, 1).Value = "Entrainement"
$ws.Cells.Item(227, 2).Value = 45874
$ws.Cells.Item(227, 3).Value = "Global"
$ws.Cells.Item(227, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(227, 6).Value = "center midfield"
$ws.Cells.Item(227, 7).Value = "01:27:50"
$ws.Cells.Item(227, 8).Value = 7.59
$ws.Cells.Item(227, 9).Value = 1.39
$ws.Cells.Item(227, 10).Value = 6.2
$ws.Cells.Item(227, 11).Value = 1.22
$ws.Cells.Item(227, 12).Value = 0.17
$ws.Cells.Item(227, 13).Value = 0
$ws.Cells.Item(227, 14).Value = 0
$ws.Cells.Item(227, 15).Value = 0
$ws.Cells.Item(227, 16).Value = 5.13
$ws.Cells.Item(227, 17).Value = 24.1
$ws.Cells.Item(227, 18).Value = 4.26
$ws.Cells.Item(227, 19).Value = 30
$ws.Cells.Item(227, 20).Value = 3
$ws.Cells.Item(227, 21).Value = 12
$ws.Cells.Item(227, 22).Value = 2

# Row 228
$ws.Cells.Item(228, 1).Value = "Entrainement"
$ws.Cells.Item(228, 2).Value = 45874
$ws.Cells.Item(228, 3).Value = "Global"
$ws.Cells.Item(228, 5).Value = "Karahali Souaré"
$ws.Cells.Item(228, 6).Value = "right forward"
$ws.Cells.Item(228, 7).Value = "01:20:18"
$ws.Cells.Item(228, 8).Value = 5.97
$ws.Cells.Item(228, 9).Value = 1.03
$ws.Cells.Item(228, 10).Value = 4.93
$ws.Cells.Item(228, 11).Value = 0.89
$ws.Cells.Item(228, 12).Value = 0.12
$ws.Cells.Item(228, 13).Value = 0.03
$ws.Cells.Item(228, 14).Value = 0
$ws.Cells.Item(228, 15).Value = 3
$ws.Cells.Item(228, 16).Value = 4.39
$ws.Cells.Item(228, 17).Value = 27.76
$ws.Cells.Item(228, 18).Value = 4.89
$ws.Cells.Item(228, 19).Value = 24
$ws.Cells.Item(228, 20).Value = 10
$ws.Cells.Item(228, 21).Value = 23
$ws.Cells.Item(228, 22).Value = 8

# Row 229
$ws.Cells.Item(229, 1).Value = "Entrainement"
$ws.Cells.Item(229, 2).Value = 45874
$ws.Cells.Item(229, 3).Value = "Global"
$ws.Cells.Item(229, 5).Value = "Amine Taiar"
$ws.Cells.Item(229, 6).Value = "center back"
$ws.Cells.Item(229, 7).Value = "01:26:37"
$ws.Cells.Item(229, 8).Value = 7.13
$ws.Cells.Item(229, 9).Value = 1.44
$ws.Cells.Item(229, 10).Value = 5.69
$ws.Cells.Item(229, 11).Value = 1.24
$ws.Cells.Item(229, 12).Value = 0.2
$ws.Cells.Item(229, 13).Value = 0
$ws.Cells.Item(229, 14).Value = 0
$ws.Cells.Item(229, 15).Value = 0
$ws.Cells.Item(229, 16).Value = 4.88
$ws.Cells.Item(229, 17).Value = 21.56
$ws.Cells.Item(229, 18).Value = 4.21
$ws.Cells.Item(229, 19).Value = 32
$ws.Cells.Item(229, 20).Value = 3
$ws.Cells.Item(229, 21).Value = 9
$ws.Cells.Item(229, 22).Value = 1

# Row 230
$ws.Cells.Item(230, 1).Value = "Entrainement"
$ws.Cells.Item(230, 2).Value = 45874
$ws.Cells.Item(230, 3).Value = "Global"
$ws.Cells.Item(230, 5).Value = "Omar Benyounes"
$ws.Cells.Item(230, 6).Value = "center midfield"
$ws.Cells.Item(230, 7).Value = "01:28:09"
$ws.Cells.Item(230, 8).Value = 6.85
$ws.Cells.Item(230, 9).Value = 1.44
$ws.Cells.Item(230, 10).Value = 5.4
$ws.Cells.Item(230, 11).Value = 1.28
$ws.Cells.Item(230, 12).Value = 0.17
$ws.Cells.Item(230, 13).Value = 0
$ws.Cells.Item(230, 14).Value = 0
$ws.Cells.Item(230, 15).Value = 0
$ws.Cells.Item(230, 16).Value = 4.45
$ws.Cells.Item(230, 17).Value = 22.57
$ws.Cells.Item(230, 18).Value = 4.24
$ws.Cells.Item(230, 19).Value = 25
$ws.Cells.Item(230, 20).Value = 1
$ws.Cells.Item(230, 21).Value = 18
$ws.Cells.Item(230, 22).Value = 5

# Update selection to mirror the final cursor position after data entry
$ws.Range("D236").Select() | Out-Null